$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B column numeric / text values that changed
$ws.Range("B1").Value = 0.0309358912037037
$ws.Range("B2").Value = 29.75749527777777
$ws.Range("B3").Value = 1501.156768779445
$ws.Range("B5").Value = 9.025
$ws.Range("B6").Value = 98
$ws.Range("B7").Value = 22

# Row 8 label + value
$ws.Range("A8").Value = "Total distance covered (km)"
$ws.Range("B8").Value = 30.62844142383255

# Row 9 label + value
$ws.Range("A9").Value = "Total energy consumption(WH/KM)"
$ws.Range("B9").Value = 49.01185626805572

# Row 10 label
$ws.Range("A10").Value = "Total SOC consumed(%)"

# Row 11 Mode text
$ws.Range("B11").Value = "Custom mode`n84.72%`nEco mode`n14.25%`nSports mode`n0.04%"

# Row 12 label
$ws.Range("A12").Value = "Peak Power(kW)"

# Row 13 label + value
$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("B13").Value = -2028.590228080331

# Row 14 label
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"

# Row 15 label + value
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 0.0008547631976997756

# Rows 16/17 swap label+value (Lowest/Highest Cell Voltage)
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.329
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 3.017

# Row 18 label
$ws.Range("A18").Value = "Difference in Cell Voltage(V)"

# Row 19 label
$ws.Range("A19").Value = "Minimum Temperature(C)"

# Row 20 label
$ws.Range("A20").Value = "Maximum Temperature(C)"

# Row 21 label + value
$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 14

# Row 22 label
$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"

# Row 23 label
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"

# Row 24 label
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"

# Row 25 label
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"

# Row 26 label
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"

# Row 27 label
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"

# Rows 28/29 swap label (lowest/highest cell temp), values stay in place
$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("A29").Value = "lowest cell temp(C)"

# Row 30 label
$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

# Row 31 previously "Maximum BMS Temperature in C" (62) is removed; data shifts up one row.
# Build the shifted block rows 31-42 (old rows 32-42 content, relabeled), then add new row 43.
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 53

$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.577147249722222

$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = [double]"1.63958255335394e-07"

$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 7.848300255452938

$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 10.55610139516604

$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 5.380231872666536

$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 8.630379249361368

$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 13.29534289644331

$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 12.1988602868933

$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 12.58793476124976

$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 19.68952642955394

$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 9.660051090587542

# New row 43
$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
